# S22/G01 - Holdings screener & batch actions
# Adds sprint S22 rows (G01 Holdings screener, G02 Groups/baskets) to the sprint task tracker,
# restyles row 179 to match the rest of the table, and adjusts a few row heights.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row height tweaks for existing rows 174-178 (minor re-wrap adjustments) ---
$ws.Rows.Item(174).RowHeight = 41.25
$ws.Rows.Item(175).RowHeight = 54.75
$ws.Rows.Item(176).RowHeight = 41.25
$ws.Rows.Item(177).RowHeight = 54.75
$ws.Rows.Item(178).RowHeight = 41.25

# --- Bring row 179 formatting in line with the rest of the table (wrap/top align, row height) ---
$ws.Range("A179:H179").WrapText = $true
$ws.Range("A179:H179").VerticalAlignment = -4160
$ws.Rows.Item(179).RowHeight = 41.75

# --- New sprint S22 rows: Holdings screener & batch actions (G01) + Groups/baskets (G02) ---

# Row 179: S21_G03_TD001
$ws.Cells.Item(179, 1).Value = "S21"
$ws.Cells.Item(179, 2).Value = "G03"
$ws.Cells.Item(179, 3).Value = "Portfolio improvement guidelines & triage framework"
$ws.Cells.Item(179, 4).Value = "S21_G03_TD001"
$ws.Cells.Item(179, 5).Value = "Document practical portfolio-stabilisation and profit-framework suggestions (A/B/C buckets, bracket use, risk sizing) in pf_improvement_suggestions.md."
$ws.Cells.Item(179, 6).Value = "Pure documentation/design task; no code changes beyond the new markdown file."
$ws.Cells.Item(179, 7).Value = "implemented"
$ws.Cells.Item(179, 8).Value = "Provides a reference playbook for using SigmaTrader tools to manage existing losers and structure new trades systematically."

# Row 180: S22_G01_TF001
$ws.Cells.Item(180, 1).Value = "S22"
$ws.Cells.Item(180, 2).Value = "G01"
$ws.Cells.Item(180, 3).Value = "Holdings screener & batch actions"
$ws.Cells.Item(180, 4).Value = "S22_G01_TF001"
$ws.Cells.Item(180, 5).Value = "Replace Holdings advanced filters with a Screener panel (builder mode) that supports multiple conditions with AND/OR and persists as a saved screener definition."
$ws.Cells.Item(180, 7).Value = "implemented"
$ws.Cells.Item(180, 8).Value = "First iteration uses client-side evaluation on the existing holdings payload and reuses DataGrid filtering."
$ws.Range("A180:E180").WrapText = $true
$ws.Range("A180:E180").VerticalAlignment = -4160
$ws.Range("G180:H180").WrapText = $true
$ws.Range("G180:H180").VerticalAlignment = -4160
$ws.Rows.Item(180).RowHeight = 41.75

# Row 181: S22_G01_TF002
$ws.Cells.Item(181, 1).Value = "S22"
$ws.Cells.Item(181, 2).Value = "G01"
$ws.Cells.Item(181, 3).Value = "Holdings screener & batch actions"
$ws.Cells.Item(181, 4).Value = "S22_G01_TF002"
$ws.Cells.Item(181, 5).Value = "Add a DSL mode to the Holdings screener that lets users write filter expressions using the indicator/alert DSL and validates them via the existing DSL parser."
$ws.Cells.Item(181, 7).Value = "implemented"
$ws.Cells.Item(181, 8).Value = "Keeps DSL syntax consistent with indicator alerts while limiting evaluation to metrics already computed for holdings."
$ws.Range("A181:E181").WrapText = $true
$ws.Range("A181:E181").VerticalAlignment = -4160
$ws.Range("G181:H181").WrapText = $true
$ws.Range("G181:H181").VerticalAlignment = -4160
$ws.Rows.Item(181).RowHeight = 41.75

# Row 182: S22_G01_TB001
$ws.Cells.Item(182, 1).Value = "S22"
$ws.Cells.Item(182, 2).Value = "G01"
$ws.Cells.Item(182, 3).Value = "Holdings screener & batch actions"
$ws.Cells.Item(182, 4).Value = "S22_G01_TB001"
$ws.Cells.Item(182, 5).Value = "Expose a backend helper or endpoint to evaluate DSL-based screener expressions against the user’s holdings, reusing the alert expression engine and indicator metrics."
$ws.Cells.Item(182, 7).Value = "implemented"
$ws.Cells.Item(182, 8).Value = "Provides a single place to evaluate more complex screeners and keeps client logic simple."
$ws.Range("A182:E182").WrapText = $true
$ws.Range("A182:E182").VerticalAlignment = -4160
$ws.Range("G182:H182").WrapText = $true
$ws.Range("G182:H182").VerticalAlignment = -4160
$ws.Rows.Item(182).RowHeight = 55.2

# Row 183: S22_G02_TB001
$ws.Cells.Item(183, 1).Value = "S22"
$ws.Cells.Item(183, 2).Value = "G02"
$ws.Cells.Item(183, 3).Value = "Groups/baskets & basic watchlists (Phase 1)"
$ws.Cells.Item(183, 4).Value = "S22_G02_TB001"
$ws.Cells.Item(183, 5).Value = "Add schema and ORM models for groups and group_members to represent named baskets, watchlists, and model portfolios."
$ws.Cells.Item(183, 7).Value = "planned"
$ws.Cells.Item(183, 8).Value = "Sets up data structures to attach symbols and optional target weights to user-defined groups."
$ws.Range("A183:E183").WrapText = $true
$ws.Range("A183:E183").VerticalAlignment = -4160
$ws.Range("G183:H183").WrapText = $true
$ws.Range("G183:H183").VerticalAlignment = -4160
$ws.Rows.Item(183).RowHeight = 41.75

# Row 184: S22_G02_TF001
$ws.Cells.Item(184, 1).Value = "S22"
$ws.Cells.Item(184, 2).Value = "G02"
$ws.Cells.Item(184, 3).Value = "Groups/baskets & basic watchlists (Phase 1)"
$ws.Cells.Item(184, 4).Value = "S22_G02_TF001"
$ws.Cells.Item(184, 5).Value = "Add a Groups page to create and edit groups, manage member symbols, and assign equal or custom target weights."
$ws.Cells.Item(184, 7).Value = "planned"
$ws.Cells.Item(184, 8).Value = "Provides a base UI for constructing candidate portfolios and watchlists independent of current holdings."
$ws.Range("A184:E184").WrapText = $true
$ws.Range("A184:E184").VerticalAlignment = -4160
$ws.Range("G184:H184").WrapText = $true
$ws.Range("G184:H184").VerticalAlignment = -4160
$ws.Rows.Item(184).RowHeight = 41.75

# Row 185: S22_G02_TF002
$ws.Cells.Item(185, 1).Value = "S22"
$ws.Cells.Item(185, 2).Value = "G02"
$ws.Cells.Item(185, 3).Value = "Groups/baskets & basic watchlists (Phase 1)"
$ws.Cells.Item(185, 4).Value = "S22_G02_TF002"
$ws.Cells.Item(185, 5).Value = "Integrate group membership into Holdings so rows can be tagged with groups and filtered by group, and add an action to allocate a fixed amount equally across a selected group via queued orders."
$ws.Cells.Item(185, 7).Value = "planned"
$ws.Cells.Item(185, 8).Value = "Reuses the existing Buy/Sell dialog and manual queue to implement basket-level investment flows."
$ws.Range("A185:E185").WrapText = $true
$ws.Range("A185:E185").VerticalAlignment = -4160
$ws.Range("G185:H185").WrapText = $true
$ws.Range("G185:H185").VerticalAlignment = -4160
$ws.Rows.Item(185).RowHeight = 55.2

# --- Update view: scroll position + active cell to match the edited area ---
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 173
$ws.Range("E179").Select()

Write-Output "S22/G01 rows added"
